$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1216
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 5041
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("F11").Value = 680
$ws.Range("F12").Value = 76

# Sheet "演出" (sheet2)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0

# Sheet "全部类型" (sheet4)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1216
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 352
$ws.Range("F7").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 9528
$ws.Range("F13").Value = 89
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = 680
